# Burndown chart update: add "Sprint 6" row of data, update existing sprint
# values, add running-total formulas, extend the chart series to include
# the new row, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing sprint values (Actual / Ideal tasks remaining) ---
$ws.Range("D5").Value = 66
$ws.Range("E5").Value = 66

$ws.Range("D6").Value = 60
$ws.Range("E6").Value = 55

$ws.Range("D7").Value = 44
$ws.Range("E7").Value = 44

$ws.Range("D8").Value = 32
$ws.Range("E8").Value = 33

$ws.Range("D9").Value = 18
$ws.Range("E9").Value = 22

# --- Add the new "Sprint 6" row ---
$ws.Range("C10").Value = "Sprint 6"
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 11

# --- Add running-total formulas under the data ---
$ws.Range("D12").Formula = "=SUM(D5:D11)"
$ws.Range("E12").Formula = "=SUM(E5:E11)"

# --- Extend the burndown chart series to cover the new Sprint 6 row ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(Sheet1!`$D`$4,Sheet1!`$C`$5:`$C`$10,Sheet1!`$D`$5:`$D`$10,1)"

$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(Sheet1!`$E`$4,,Sheet1!`$E`$5:`$E`$10,2)"

# --- Move the active selection, matching the authored workbook state ---
$ws.Range("S35").Select()
